{"js": "// Merge the first \"Section Headnote\" paragraph (\"This is the first chapter\n// of the casebook.\") together with every paragraph through the second\n// \"Section Headnote\" paragraph (\"This is the second chapter of the\n// casebook.\") into a single paragraph/run. The \"headnote\" commentary\n// paragraphs' own text (ResourceHeadnote / SectionHeadnote) is dropped;\n// only the surrounding number/title/body text survives, concatenated with\n// no separators, and the trailing \"Headnote\" paragraph is removed outright.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Locate the anchor paragraphs by their distinctive text rather than a\n// hard-coded index, so the script is resilient to minor doc changes.\nlet startIndex = -1;\nlet endIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (startIndex === -1 && t.indexOf(\"This is the first chapter of the casebook.\") !== -1) {\n    startIndex = i;\n  }\n  if (t.indexOf(\"This is the second chapter of the casebook.\") !== -1) {\n    endIndex = i;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1 || endIndex <= startIndex) {\n  throw new Error(\"Could not locate the paragraph range to merge.\");\n}\n\nconst startStyle = paragraphs.items[startIndex].style;\n\n// Concatenate the text of every paragraph strictly between the two anchors\n// (the anchors themselves are handled separately below), skipping\n// \"Headnote\" commentary paragraphs whose own text is dropped entirely.\nlet mergedText = \"\";\nfor (let i = startIndex + 1; i < endIndex; i++) {\n  const style = paragraphs.items[i].style;\n  if (style && style.indexOf(\"Headnote\") !== -1) {\n    continue;\n  }\n  mergedText += paragraphs.items[i].text;\n}\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Rewrite the first anchor paragraph via OOXML injection (rather than\n// insertText) so that the embedded \"\\n\" character lands as a literal\n// character inside <w:t>, not as a new-paragraph break.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"' +\n  startStyle.replace(/\\s+/g, \"\") +\n  '\"/></w:pPr><w:r><w:t xml:space=\"preserve\">' +\n  xmlEscape(mergedText) +\n  \"</w:t></w:r></w:p></w:body></w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst startRange = paragraphs.items[startIndex].getRange(\"Whole\");\nstartRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the paragraphs (the collection is now stale after the OOXML\n// insert) and delete everything from just after the rewritten paragraph\n// through the second anchor (inclusive), walking backwards so indices\n// stay valid as each delete happens.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet newStartIndex = -1;\nlet newEndIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const t = paragraphs2.items[i].text;\n  if (newStartIndex === -1 && t.indexOf(mergedText) !== -1) {\n    newStartIndex = i;\n  }\n  if (t.indexOf(\"This is the second chapter of the casebook.\") !== -1) {\n    newEndIndex = i;\n  }\n}\n\nfor (let i = newEndIndex; i > newStartIndex; i--) {\n  paragraphs2.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Merge the first \"Section Headnote\" paragraph (\"This is the first chapter\n# of the casebook.\") together with every paragraph through the second\n# \"Section Headnote\" paragraph (\"This is the second chapter of the\n# casebook.\") into a single paragraph/run. The \"headnote\" commentary\n# paragraphs' own text (Resource Headnote / Section Headnote) is dropped;\n# only the surrounding number/title/body text survives, concatenated with\n# no separators, and the trailing headnote paragraph is removed outright.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$startIdx = -1\n$endIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($startIdx -eq -1 -and $t -like \"*This is the first chapter of the casebook.*\") {\n        $startIdx = $i\n    }\n    if ($t -like \"*This is the second chapter of the casebook.*\") {\n        $endIdx = $i\n    }\n}\n\nif ($startIdx -eq -1 -or $endIdx -eq -1 -or $endIdx -le $startIdx) {\n    throw \"Could not locate the paragraph range to merge.\"\n}\n\n# Concatenate the text of every paragraph strictly between the two anchors,\n# skipping \"Headnote\" commentary paragraphs whose own text is dropped\n# entirely. Paragraph.Range.Text always carries a trailing CR (U+000D)\n# paragraph-mark character that is not part of the visible text, so trim\n# exactly one off the end of each paragraph's text before concatenating --\n# any genuine embedded \"\\n\" inside the paragraph's own text is left alone.\n$merged = \"\"\nfor ($i = $startIdx + 1; $i -lt $endIdx; $i++) {\n    $p = $d.Paragraphs($i)\n    $styleName = $p.Style.NameLocal\n    if ($styleName -like \"*Headnote*\") {\n        continue\n    }\n    $merged += $p.Range.Text.TrimEnd([char]13)\n}\n\n# Delete every paragraph from just after the first anchor through the\n# second anchor (inclusive) as a single range delete, done before the text\n# assignment so the first anchor paragraph's index/range stays valid.\n$deleteStart = $d.Paragraphs($startIdx + 1).Range.Start\n$deleteEnd = $d.Paragraphs($endIdx).Range.End\n$d.Range($deleteStart, $deleteEnd).Delete()\n\n# Replace the first anchor paragraph's text with the merged content.\n$d.Paragraphs($startIdx).Range.Text = $merged\n"}
